$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value2 = '20.601.57'
$ws.Range('E2').Value2 = '  +0.41%  '
$ws.Range('D3').Value2 = '1.481.73'
$ws.Range('E3').Value2 = '  +0.69%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '1.011'
$ws.Range('E4').Value2 = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '0.9711'
$ws.Range('E5').Value2 = '  +2.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '280.61'
$ws.Range('E6').Value2 = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.3665'
$ws.Range('E7').Value2 = '  -1.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.3088'
$ws.Range('E8').Value2 = '  -3.34%  '
$ws.Range('E9').Value2 = '  -3.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '1.064'
$ws.Range('E10').Value2 = '  +0.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '0.06679'
$ws.Range('E11').Value2 = '  +0.06%  '
$ws.Range('E12').Value2 = '  +0.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '5.540'
$ws.Range('E13').Value2 = '  -1.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '18.18'
$ws.Range('E14').Value2 = '  -0.62%  '
$ws.Range('E15').Value2 = '  -0.16%  '
$ws.Range('E16').Value2 = '  +2.37%  '
$ws.Range('E17').Value2 = '  -0.20%  '
$ws.Range('D18').Value2 = '1.483.87'
$ws.Range('E18').Value2 = '  +0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '0.05962'
$ws.Range('E19').Value2 = '  +3.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '70.01'
$ws.Range('E20').Value2 = '  -3.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '5.513'
$ws.Range('E21').Value2 = '  -3.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '14.56'
$ws.Range('E22').Value2 = '  -1.46%  '
$ws.Range('E23').Value2 = '  -1.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '2.275'
$ws.Range('E24').Value2 = '  +0.03%  '
$ws.Range('D25').Value2 = '20.637.52'
$ws.Range('E25').Value2 = '  -0.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '142.49'
$ws.Range('E26').Value2 = '  +3.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '2.126'
$ws.Range('E27').Value2 = '  -7.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '17.34'
$ws.Range('E28').Value2 = '  -1.57%  '
$ws.Range('D29').Value2 = '1.646.30'
$ws.Range('E29').Value2 = '  +0.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '114.35'
$ws.Range('E30').Value2 = '  +0.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '3.966'
$ws.Range('E31').Value2 = '  +0.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '0.8306'
$ws.Range('E32').Value2 = '  -2.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '5.046'
$ws.Range('E33').Value2 = '  -5.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '0.08032'
$ws.Range('E34').Value2 = '  +2.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '1.546'
$ws.Range('E35').Value2 = '  -4.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '1.220'
$ws.Range('E36').Value2 = '  +8.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '0.05811'
$ws.Range('E37').Value2 = '  -4.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '4.771'
$ws.Range('E38').Value2 = '  -3.49%  '
$ws.Range('B39').Value2 = 'VeChain'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '0.02053'
$ws.Range('E39').Value2 = '  -0.99%  '
$ws.Range('B40').Value2 = 'Frax'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '0.9706'
$ws.Range('E40').Value2 = '  +1.37%  '
$ws.Range('B41').Value2 = 'FraxShare'
$ws.Range('C41').Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '7.695'
$ws.Range('E41').Value2 = '  +1.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '10.46'
$ws.Range('E42').Value2 = '  -2.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '0.1886'
$ws.Range('E43').Value2 = '  -1.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '0.5324'
$ws.Range('E44').Value2 = '  -1.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '3.547'
$ws.Range('E45').Value2 = '  -1.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '12.25'
$ws.Range('E46').Value2 = '  -2.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '119.32'
$ws.Range('E47').Value2 = '  -2.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '0.5220'
$ws.Range('E48').Value2 = '  -2.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '1.838'
$ws.Range('E49').Value2 = '  +0.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '0.06516'
$ws.Range('E50').Value2 = '  +0.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '0.9922'
$ws.Range('E51').Value2 = '  -0.02%  '
